$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

for ($row = 6; $row -le 34; $row++) {
    $cell = $ws.Range("AH$row")
    if ($cell.Value2 -eq "Cellranger v3.0.2") {
        $cell.Value2 = "Cellranger v3.0.1"
    }
}
